# Apply the diff: add three new log rows (13-15) for LeechT neuron runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "2018.03.08" looks like a date to Excel's smart-parser, which would turn it
# into a date serial number instead of leaving it as plain text (which is how
# the existing "Date" column values are stored in the workbook). Build it via
# a formula that evaluates to text, then copy/paste-special (values only)
# directly into each destination cell so it ends up holding a genuine text
# string without picking up any stray number-format/style.
$tmp = $ws.Cells.Item(200, 200)
$tmp.Formula = "=""2018.03.08"""
$tmp.Copy() | Out-Null
$ws.Cells.Item(13, 1).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(14, 1).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(15, 1).PasteSpecial(-4163) | Out-Null
$tmp.Clear()

$rows = @(
    @{ Row = 13; Time = "17:47:10"; Neuron = "LeechT"; Tstim = 20; PRF = 310; DutyFactor = "N/A"; Astim = 0.8;  Spikes = 12; Latency = 3.501129396579542; SpikeRate = 0.03874718679669917 },
    @{ Row = 14; Time = "17:47:36"; Neuron = "LeechT"; Tstim = 20; PRF = 310; DutyFactor = "N/A"; Astim = 0.82; Spikes = 12; Latency = 3.501129396579542; SpikeRate = 0.03874718679669917 },
    @{ Row = 15; Time = "17:47:43"; Neuron = "LeechT"; Tstim = 20; PRF = 300; DutyFactor = "N/A"; Astim = 0.8;  Spikes = 12; Latency = 3.501167055685228; SpikeRate = 0.0387467700258398 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.Time
    $ws.Cells.Item($row, 3).Value = $r.Neuron
    $ws.Cells.Item($row, 4).Value = $r.Tstim
    $ws.Cells.Item($row, 5).Value = $r.PRF
    $ws.Cells.Item($row, 6).Value = $r.DutyFactor
    $ws.Cells.Item($row, 7).Value = 1
    $ws.Cells.Item($row, 8).Value = 3499
    $ws.Cells.Item($row, 9).Value = $r.Astim
    $ws.Cells.Item($row, 10).Value = $r.Spikes
    $ws.Cells.Item($row, 11).Value = $r.Latency
    $ws.Cells.Item($row, 12).Value = $r.SpikeRate
}
